# Apply the claim-text revisions described by the diff.
$d = $word.ActiveDocument

# 1) First paragraph: collapse the "방법으로서, 다음의 단계를 포함하는 방법:" framing
#    into the shorter "특성화하는 방법에 있어서," framing.
$d.Content.Find.Execute(
    "단백질의 특성을 규명하는 방법으로서, 다음의 단계를 포함하는 방법:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "단백질을 특성화하는 방법에 있어서,", 2)

# 2) Reorder "상기 시료를 분광법을 위해" -> "분광법을 위해 상기 시료를"
$d.Content.Find.Execute(
    "상기 시료를 분광법을 위해 준비하는 단계;",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "분광법을 위해 상기 시료를 준비하는 단계;", 2)

# 3) "빈 영역에서" -> "빈 영역으로부터"
$d.Content.Find.Execute(
    "결과로 얻어진 스펙트럼의 빈 영역에서 잡음을 제거하는 단계; 및",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "결과로 얻어진 스펙트럼의 빈 영역으로부터 잡음을 제거하는 단계; 및", 2)

# 4) Rephrase the final step sentence.
$d.Content.Find.Execute(
    "상기 스펙트럼을 분석하여 단백질의 특성을 규명하는 단계.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "상기 단백질을 특성화하기 위해 상기 스펙트럼을 분석하는 단계를 포함하는 방법.", 2)
